$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 197.5433703333333
$ws.Range("H2").Value = 592.6301109999999
$ws.Range("I2").Value = 0.3388703761585983
$ws.Range("J2").Value = 0.3388703761585982
$ws.Range("M2").Value = 2.385742333333333
$ws.Range("N2").Value = 7.157227
$ws.Range("O2").Value = 0.2243329847197944
$ws.Range("P2").Value = 0.2243329847197944
$ws.Range("Q2").Value = 471.2875812735774
$ws.Range("R2").Value = 4241.588231462197
$ws.Range("S2").Value = 0.0760198029167778
$ws.Range("T2").Value = 0.07601980291677779
# Row 3
$ws.Range("G3").Value = 197.5433703333333
$ws.Range("H3").Value = 592.6301109999999
$ws.Range("I3").Value = 0.3388703761585983
$ws.Range("J3").Value = 0.3388703761585982
$ws.Range("O3").Value = 0.01121613676875902
$ws.Range("P3").Value = 0.01121613676875902
$ws.Range("Q3").Value = 23.56330245231055
$ws.Range("R3").Value = 212.069722070795
$ws.Range("S3").Value = 0.003800816485875654
$ws.Range("T3").Value = 0.003800816485875653
# Row 4
$ws.Range("G4").Value = 197.5433703333333
$ws.Range("H4").Value = 592.6301109999999
$ws.Range("I4").Value = 0.3388703761585983
$ws.Range("J4").Value = 0.3388703761585982
$ws.Range("M4").Value = 5.307525999999999
$ws.Range("N4").Value = 15.922578
$ws.Range("O4").Value = 0.4990703029502535
$ws.Range("P4").Value = 0.4990703029502535
$ws.Range("Q4").Value = 1048.466574171795
$ws.Range("R4").Value = 9436.199167546156
$ws.Range("S4").Value = 0.169120141290338
$ws.Range("T4").Value = 0.1691201412903379
# Row 5
$ws.Range("G5").Value = 197.5433703333333
$ws.Range("H5").Value = 592.6301109999999
$ws.Range("I5").Value = 0.3388703761585983
$ws.Range("J5").Value = 0.3388703761585982
$ws.Range("M5").Value = 0.189174
$ws.Range("N5").Value = 0.5675220000000001
$ws.Range("O5").Value = 0.01778816071561614
$ws.Range("P5").Value = 0.01778816071561614
$ws.Range("Q5").Value = 37.37006953943801
$ws.Range("R5").Value = 336.330625854942
$ws.Range("S5").Value = 0.006027880712870442
$ws.Range("T5").Value = 0.006027880712870441
# Row 6
$ws.Range("G6").Value = 197.5433703333333
$ws.Range("H6").Value = 592.6301109999999
$ws.Range("I6").Value = 0.3388703761585983
$ws.Range("J6").Value = 0.3388703761585982
$ws.Range("M6").Value = 2.633102333333333
$ws.Range("N6").Value = 7.899307
$ws.Range("O6").Value = 0.247592414845577
$ws.Range("P6").Value = 0.247592414845577
$ws.Range("Q6").Value = 520.1519093592308
$ws.Range("R6").Value = 4681.367184233077
$ws.Range("S6").Value = 0.08390173475273641
$ws.Range("T6").Value = 0.08390173475273639
# Row 7
$ws.Range("I7").Value = 0.1369374790620155
$ws.Range("J7").Value = 0.1369374790620154
$ws.Range("M7").Value = 2.385742333333333
$ws.Range("N7").Value = 7.157227
$ws.Range("O7").Value = 0.2243329847197944
$ws.Range("P7").Value = 0.2243329847197944
$ws.Range("Q7").Value = 190.4472560405629
$ws.Range("R7").Value = 1714.025304365066
$ws.Range("S7").Value = 0.03071959339798628
$ws.Range("T7").Value = 0.03071959339798628
# Row 8
$ws.Range("I8").Value = 0.1369374790620155
$ws.Range("J8").Value = 0.1369374790620154
$ws.Range("O8").Value = 0.01121613676875902
$ws.Range("P8").Value = 0.01121613676875902
$ws.Range("S8").Value = 0.00153590949392864
$ws.Range("T8").Value = 0.00153590949392864
# Row 9
$ws.Range("I9").Value = 0.1369374790620155
$ws.Range("J9").Value = 0.1369374790620154
$ws.Range("M9").Value = 5.307525999999999
$ws.Range("N9").Value = 15.922578
$ws.Range("O9").Value = 0.4990703029502535
$ws.Range("P9").Value = 0.4990703029502535
$ws.Range("Q9").Value = 423.6852190369026
$ws.Range("R9").Value = 3813.166971332123
$ws.Range("S9").Value = 0.06834142916072405
$ws.Range("T9").Value = 0.06834142916072404
# Row 10
$ws.Range("I10").Value = 0.1369374790620155
$ws.Range("J10").Value = 0.1369374790620154
$ws.Range("M10").Value = 0.189174
$ws.Range("N10").Value = 0.5675220000000001
$ws.Range("O10").Value = 0.01778816071561614
$ws.Range("P10").Value = 0.01778816071561614
$ws.Range("Q10").Value = 15.101240695964
$ws.Range("R10").Value = 135.911166263676
$ws.Range("S10").Value = 0.002435865885546451
$ws.Range("T10").Value = 0.002435865885546451
# Row 11
$ws.Range("I11").Value = 0.1369374790620155
$ws.Range("J11").Value = 0.1369374790620154
$ws.Range("M11").Value = 2.633102333333333
$ws.Range("N11").Value = 7.899307
$ws.Range("O11").Value = 0.247592414845577
$ws.Range("P11").Value = 0.247592414845577
$ws.Range("Q11").Value = 210.1933252601896
$ws.Range("R11").Value = 1891.739927341706
$ws.Range("S11").Value = 0.03390468112383006
$ws.Range("T11").Value = 0.03390468112383005
# Row 12
$ws.Range("G12").Value = 148.824417
$ws.Range("H12").Value = 446.473251
$ws.Range("I12").Value = 0.2552967790580629
$ws.Range("J12").Value = 0.2552967790580629
$ws.Range("M12").Value = 2.385742333333333
$ws.Range("N12").Value = 7.157227
$ws.Range("O12").Value = 0.2243329847197944
$ws.Range("P12").Value = 0.2243329847197944
$ws.Range("Q12").Value = 355.056711870553
$ws.Range("R12").Value = 3195.510406834977
$ws.Range("S12").Value = 0.05727148843544513
$ws.Range("T12").Value = 0.05727148843544513
# Row 13
$ws.Range("G13").Value = 148.824417
$ws.Range("H13").Value = 446.473251
$ws.Range("I13").Value = 0.2552967790580629
$ws.Range("J13").Value = 0.2552967790580629
$ws.Range("O13").Value = 0.01121613676875902
$ws.Range("P13").Value = 0.01121613676875902
$ws.Range("Q13").Value = 17.752024500455
$ws.Range("R13").Value = 159.768220504095
$ws.Range("S13").Value = 0.002863443590538887
$ws.Range("T13").Value = 0.002863443590538887
# Row 14
$ws.Range("G14").Value = 148.824417
$ws.Range("H14").Value = 446.473251
$ws.Range("I14").Value = 0.2552967790580629
$ws.Range("J14").Value = 0.2552967790580629
$ws.Range("M14").Value = 5.307525999999999
$ws.Range("N14").Value = 15.922578
$ws.Range("O14").Value = 0.4990703029502535
$ws.Range("P14").Value = 0.4990703029502535
$ws.Range("Q14").Value = 789.889462662342
$ws.Range("R14").Value = 7109.005163961077
$ws.Range("S14").Value = 0.1274110408667314
$ws.Range("T14").Value = 0.1274110408667314
# Row 15
$ws.Range("G15").Value = 148.824417
$ws.Range("H15").Value = 446.473251
$ws.Range("I15").Value = 0.2552967790580629
$ws.Range("J15").Value = 0.2552967790580629
$ws.Range("M15").Value = 0.189174
$ws.Range("N15").Value = 0.5675220000000001
$ws.Range("O15").Value = 0.01778816071561614
$ws.Range("P15").Value = 0.01778816071561614
$ws.Range("Q15").Value = 28.15371026155801
$ws.Range("R15").Value = 253.383392354022
$ws.Range("S15").Value = 0.004541260136063967
$ws.Range("T15").Value = 0.004541260136063967
# Row 16
$ws.Range("G16").Value = 148.824417
$ws.Range("H16").Value = 446.473251
$ws.Range("I16").Value = 0.2552967790580629
$ws.Range("J16").Value = 0.2552967790580629
$ws.Range("M16").Value = 2.633102333333333
$ws.Range("N16").Value = 7.899307
$ws.Range("O16").Value = 0.247592414845577
$ws.Range("P16").Value = 0.247592414845577
$ws.Range("Q16").Value = 391.869919659673
$ws.Range("R16").Value = 3526.829276937057
$ws.Range("S16").Value = 0.06320954602928353
$ws.Range("T16").Value = 0.06320954602928353
# Row 17
$ws.Range("G17").Value = 35.426853
$ws.Range("H17").Value = 106.280559
$ws.Range("I17").Value = 0.06077202683121193
$ws.Range("J17").Value = 0.06077202683121192
$ws.Range("M17").Value = 2.385742333333333
$ws.Range("N17").Value = 7.157227
$ws.Range("O17").Value = 0.2243329847197944
$ws.Range("P17").Value = 0.2243329847197944
$ws.Range("Q17").Value = 84.519342938877
$ws.Range("R17").Value = 760.6740864498929
$ws.Range("S17").Value = 0.0136331701665172
$ws.Range("T17").Value = 0.0136331701665172
# Row 18
$ws.Range("G18").Value = 35.426853
$ws.Range("H18").Value = 106.280559
$ws.Range("I18").Value = 0.06077202683121193
$ws.Range("J18").Value = 0.06077202683121192
$ws.Range("O18").Value = 0.01121613676875902
$ws.Range("P18").Value = 0.01121613676875902
$ws.Range("Q18").Value = 4.225774070595
$ws.Range("R18").Value = 38.031966635355
$ws.Range("S18").Value = 0.0006816273646535658
$ws.Range("T18").Value = 0.0006816273646535657
# Row 19
$ws.Range("G19").Value = 35.426853
$ws.Range("H19").Value = 106.280559
$ws.Range("I19").Value = 0.06077202683121193
$ws.Range("J19").Value = 0.06077202683121192
$ws.Range("M19").Value = 5.307525999999999
$ws.Range("N19").Value = 15.922578
$ws.Range("O19").Value = 0.4990703029502535
$ws.Range("P19").Value = 0.4990703029502535
$ws.Range("Q19").Value = 188.028943395678
$ws.Range("R19").Value = 1692.260490561102
$ws.Range("S19").Value = 0.03032951384155387
$ws.Range("T19").Value = 0.03032951384155386
# Row 20
$ws.Range("G20").Value = 35.426853
$ws.Range("H20").Value = 106.280559
$ws.Range("I20").Value = 0.06077202683121193
$ws.Range("J20").Value = 0.06077202683121192
$ws.Range("M20").Value = 0.189174
$ws.Range("N20").Value = 0.5675220000000001
$ws.Range("O20").Value = 0.01778816071561614
$ws.Range("P20").Value = 0.01778816071561614
$ws.Range("Q20").Value = 6.701839489422001
$ws.Range("R20").Value = 60.31655540479801
$ws.Range("S20").Value = 0.001081022580287334
$ws.Range("T20").Value = 0.001081022580287334
# Row 21
$ws.Range("G21").Value = 35.426853
$ws.Range("H21").Value = 106.280559
$ws.Range("I21").Value = 0.06077202683121193
$ws.Range("J21").Value = 0.06077202683121192
$ws.Range("M21").Value = 2.633102333333333
$ws.Range("N21").Value = 7.899307
$ws.Range("O21").Value = 0.247592414845577
$ws.Range("P21").Value = 0.247592414845577
$ws.Range("Q21").Value = 93.28252929695701
$ws.Range("R21").Value = 839.542763672613
$ws.Range("S21").Value = 0.01504669287819996
$ws.Range("T21").Value = 0.01504669287819996
# Row 22
$ws.Range("G22").Value = 121.3248153333333
$ws.Range("H22").Value = 363.974446
$ws.Range("I22").Value = 0.2081233388901116
$ws.Range("J22").Value = 0.2081233388901115
$ws.Range("M22").Value = 2.385742333333333
$ws.Range("N22").Value = 7.157227
$ws.Range("O22").Value = 0.2243329847197944
$ws.Range("P22").Value = 0.2243329847197944
$ws.Range("Q22").Value = 289.4497480245824
$ws.Range("R22").Value = 2605.047732221242
$ws.Range("S22").Value = 0.04668892980306798
$ws.Range("T22").Value = 0.04668892980306798
# Row 23
$ws.Range("G23").Value = 121.3248153333333
$ws.Range("H23").Value = 363.974446
$ws.Range("I23").Value = 0.2081233388901116
$ws.Range("J23").Value = 0.2081233388901115
$ws.Range("O23").Value = 0.01121613676875902
$ws.Range("P23").Value = 0.01121613676875902
$ws.Range("Q23").Value = 14.47182618098556
$ws.Range("R23").Value = 130.24643562887
$ws.Range("S23").Value = 0.002334339833762274
$ws.Range("T23").Value = 0.002334339833762274
# Row 24
$ws.Range("G24").Value = 121.3248153333333
$ws.Range("H24").Value = 363.974446
$ws.Range("I24").Value = 0.2081233388901116
$ws.Range("J24").Value = 0.2081233388901115
$ws.Range("M24").Value = 5.307525999999999
$ws.Range("N24").Value = 15.922578
$ws.Range("O24").Value = 0.4990703029502535
$ws.Range("P24").Value = 0.4990703029502535
$ws.Range("Q24").Value = 643.9346118268652
$ws.Range("R24").Value = 5795.411506441787
$ws.Range("S24").Value = 0.1038681777909062
$ws.Range("T24").Value = 0.1038681777909062
# Row 25
$ws.Range("G25").Value = 121.3248153333333
$ws.Range("H25").Value = 363.974446
$ws.Range("I25").Value = 0.2081233388901116
$ws.Range("J25").Value = 0.2081233388901115
$ws.Range("M25").Value = 0.189174
$ws.Range("N25").Value = 0.5675220000000001
$ws.Range("O25").Value = 0.01778816071561614
$ws.Range("P25").Value = 0.01778816071561614
$ws.Range("Q25").Value = 22.951500615868
$ws.Range("R25").Value = 206.563505542812
$ws.Range("S25").Value = 0.003702131400847947
$ws.Range("T25").Value = 0.003702131400847947
# Row 26
$ws.Range("G26").Value = 121.3248153333333
$ws.Range("H26").Value = 363.974446
$ws.Range("I26").Value = 0.2081233388901116
$ws.Range("J26").Value = 0.2081233388901115
$ws.Range("M26").Value = 2.633102333333333
$ws.Range("N26").Value = 7.899307
$ws.Range("O26").Value = 0.247592414845577
$ws.Range("P26").Value = 0.247592414845577
$ws.Range("Q26").Value = 319.4606543454358
$ws.Range("R26").Value = 2875.145889108922
$ws.Range("S26").Value = 0.05152976006152712
$ws.Range("T26").Value = 0.05152976006152711
